# Update DRC to remove reference to EV charging
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Revised Calcs_India": core numeric/formula edits + new note rows
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Revised Calcs_India")
$wsCalc.Activate() | Out-Null

# J2 used to be a literal 90; now it nets off 12.5 (EV-charging contribution)
$wsCalc.Range("J2").Formula = "=90-12.5"

# C3 used to be a literal 90000; now it is derived straight from J2
$wsCalc.Range("C3").Formula = "=J2*1000"

# Insert two fresh rows just above the "Time (Year)" table (old row 7) so the
# existing data block shifts from rows 7-9 down to rows 9-11, and the
# trailing source-note row shifts from 27 to 29.
$wsCalc.Rows("7:8").Insert() | Out-Null

# Move the chart/picture anchored on this sheet down by the same two rows
# (it does not follow the row insertion automatically).
$picCalc = $wsCalc.Shapes.Item(1)
$origTop = $picCalc.Top
$origLeft = $picCalc.Left
$origWidth = $picCalc.Width
$origHeight = $picCalc.Height
$picCalc.Top = $origTop + 28.5
$picCalc.Left = $origLeft
$picCalc.Width = $origWidth
$picCalc.Height = $origHeight

# Populate the two newly inserted rows with the explanatory note about
# removing the EV-charging contribution, matching the italic "source note"
# styling already used by the neighbouring notes (row 4/5 and the old row 27).
$wsCalc.Range("G4").Copy() | Out-Null
$wsCalc.Range("G6:G8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$wsCalc.Range("G6").Value = "However, we adjust to remove the contribution from EV charging, which"
$wsCalc.Range("G7").Value = "is calculated separately in the EPS."

# Update the saved selection for this sheet.
$wsCalc.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "About": new hyperlink on G8 + refreshed selection
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null

$cpiUrl = "https://www.climatepolicyinitiative.org/wp-content/uploads/2020/08/CPI-India-flexibility-25-August-2020-full-report-1.pdf"
$wsAbout.Hyperlinks.Add($wsAbout.Range("G8"), $cpiUrl) | Out-Null

$wsAbout.Range("G9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "DRC-BDRC": refreshed selection/scroll position
# ---------------------------------------------------------------------------
$wsBdrc = $wb.Worksheets.Item("DRC-BDRC")
$wsBdrc.Activate() | Out-Null
$wsBdrc.Range("Q2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 13

# ---------------------------------------------------------------------------
# Sheet "DRC-PADRC": refreshed selection/scroll position
# ---------------------------------------------------------------------------
$wsPadrc = $wb.Worksheets.Item("DRC-PADRC")
$wsPadrc.Activate() | Out-Null
$wsPadrc.Range("Q2").Select() | Out-Null

# Leave the "About" sheet as the active/visible one, matching the source file.
$wsAbout.Activate() | Out-Null

$wb.Application.CalculateFull() | Out-Null
